$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.746.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.008.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.21%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.53%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.522"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.995.95"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.96%  "
$ws.Range("E10").Value = "  +7.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.30"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +12.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.455"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.38%  "
$ws.Range("E13").Value = "  +7.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.28%  "
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.501.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.004.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "59.709.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "438.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.718"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.49%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  +11.65%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  +4.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.88"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.106"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +9.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.995"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.46%  "
$ws.Range("E36").Value = "  +4.70%  "
$ws.Range("E37").Value = "  +2.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.09"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.04%  "
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("E40").Value = "  +11.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "402.47"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.56%  "
$ws.Range("E42").Value = "  +3.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.766.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.87%  "
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("E45").Value = "  +6.95%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "123.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +22.27%  "
$ws.Range("E49").Value = "  +4.98%  "
$ws.Range("E50").Value = "  +1.98%  "
$ws.Range("E51").Value = "  +3.63%  "
